$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of |S*|/n column (J) across the 10 result rows
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary rows 14-17: labels in column A, aggregate formulas in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Row height for the new summary rows
$ws.Range("A14:A17").RowHeight = 15.6

# Style B14 (bold, size 12, vertically centered) and propagate the format
# to the rest of the summary column so only one new cell style is created
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Bold the J12 average (reuses the existing header-style bold font)
$ws.Range("J12").Font.Bold = $true

# Leave the summary block selected, matching the saved view state
$ws.Range("A14:B17").Select()

# Page setup (paper size / orientation) as saved by the newer Excel version
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "done"
